{"js": "/*\n * Edit: update the worksheet date and every answer cell in the two-digit\n * division table.\n *\n * Strategy: every run of text in this document is either (a) the single\n * date paragraph above the table, or (b) the lone run inside a table\n * cell. We replace each one by its document position (paragraph / table\n * cell index) rather than by a blind whole-document text search, because\n * several new values happen to equal the OLD value of a *different*\n * cell (e.g. cell (0,0) becomes \"31\u00f75=6, 1\", which is the original text\n * of cell (12,2)). A global search/replace done in document order would\n * therefore wrongly re-match already-updated text. Scoping each search\n * to the owning paragraph/cell avoids that collision entirely, and using\n * `insertText(..., Word.InsertLocation.replace)` on the matched (narrow)\n * range preserves the existing run/paragraph formatting (fonts, size,\n * alignment) instead of resetting it the way replacing the whole\n * paragraph/cell body would.\n */\n\nconst DATE_OLD = \"2023-07-27 Thursday\";\nconst DATE_NEW = \"2023-07-28 Friday\";\n\n// row/col are 0-based table coordinates (as used by Table.getCell).\nconst CELL_EDITS = [\n  {\n    \"row\": 0,\n    \"col\": 0,\n    \"old\": \"32\u00f74=8, 0\",\n    \"new\": \"31\u00f75=6, 1\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 1,\n    \"old\": \"29\u00f74=7, 1\",\n    \"new\": \"11\u00f74=2, 3\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 2,\n    \"old\": \"68\u00f79=7, 5\",\n    \"new\": \"47\u00f78=5, 7\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 3,\n    \"old\": \"12\u00f72=6, 0\",\n    \"new\": \"48\u00f73=16, 0\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 4,\n    \"old\": \"49\u00f73=16, 1\",\n    \"new\": \"11\u00f72=5, 1\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 0,\n    \"old\": \"46\u00f79=5, 1\",\n    \"new\": \"82\u00f76=13, 4\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 1,\n    \"old\": \"55\u00f77=7, 6\",\n    \"new\": \"13\u00f78=1, 5\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 2,\n    \"old\": \"68\u00f75=13, 3\",\n    \"new\": \"48\u00f78=6, 0\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 3,\n    \"old\": \"47\u00f78=5, 7\",\n    \"new\": \"20\u00f78=2, 4\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 4,\n    \"old\": \"59\u00f76=9, 5\",\n    \"new\": \"80\u00f73=26, 2\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 0,\n    \"old\": \"69\u00f73=23, 0\",\n    \"new\": \"35\u00f78=4, 3\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 1,\n    \"old\": \"50\u00f77=7, 1\",\n    \"new\": \"14\u00f78=1, 6\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 2,\n    \"old\": \"58\u00f78=7, 2\",\n    \"new\": \"32\u00f76=5, 2\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 3,\n    \"old\": \"31\u00f79=3, 4\",\n    \"new\": \"82\u00f76=13, 4\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 4,\n    \"old\": \"21\u00f77=3, 0\",\n    \"new\": \"42\u00f79=4, 6\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 0,\n    \"old\": \"52\u00f76=8, 4\",\n    \"new\": \"84\u00f74=21, 0\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 1,\n    \"old\": \"60\u00f76=10, 0\",\n    \"new\": \"54\u00f78=6, 6\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 2,\n    \"old\": \"31\u00f75=6, 1\",\n    \"new\": \"99\u00f72=49, 1\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 3,\n    \"old\": \"80\u00f73=26, 2\",\n    \"new\": \"32\u00f73=10, 2\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 4,\n    \"old\": \"82\u00f72=41, 0\",\n    \"new\": \"98\u00f77=14, 0\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 0,\n    \"old\": \"32\u00f77=4, 4\",\n    \"new\": \"19\u00f74=4, 3\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 1,\n    \"old\": \"81\u00f72=40, 1\",\n    \"new\": \"92\u00f75=18, 2\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 2,\n    \"old\": \"45\u00f78=5, 5\",\n    \"new\": \"67\u00f75=13, 2\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 3,\n    \"old\": \"94\u00f73=31, 1\",\n    \"new\": \"75\u00f74=18, 3\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 4,\n    \"old\": \"38\u00f72=19, 0\",\n    \"new\": \"99\u00f78=12, 3\"\n  }\n];\n\n// --- Update the date paragraph (first paragraph in the body, above the table) ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\nconst dateResults = dateParagraph.search(DATE_OLD, { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length === 0) {\n  throw new Error(`Date paragraph text not found: ${DATE_OLD}`);\n}\ndateResults.items[0].insertText(DATE_NEW, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Update every answer cell in the (single) table ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const edit of CELL_EDITS) {\n  const cell = table.getCell(edit.row, edit.col);\n  const results = cell.body.search(edit.old, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Cell (${edit.row}, ${edit.col}) text not found: ${edit.old}`);\n  }\n  results.items[0].insertText(edit.new, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date (first paragraph, above the table) and every\n# answer cell in the two-digit division table.\n#\n# Each edit is applied directly to the Range of the specific paragraph /\n# table cell it belongs to (1-based indices, as Word's object model uses),\n# rather than via a document-wide Find & Replace. Two reasons:\n#   1. Several of the NEW values equal the OLD value of a *different* cell\n#      (e.g. cell R1C1 becomes \"31\u00f75=6, 1\", which is the original text of\n#      cell R13C3), so a sequential text-based replace pass risks re-matching\n#      text that an earlier step already wrote.\n#   2. In this host, Range.Find.Execute searches the whole document body\n#      instead of staying confined to the Range it was called on, so scoping\n#      a Find/Replace to an individual cell's Range does not actually avoid\n#      the collision above.\n# Assigning Range.Text directly targets only that single paragraph/cell by\n# construction (no text search involved) and keeps the existing run/\n# paragraph formatting (fonts, size, alignment), since only the text inside\n# the run is swapped, not the run/paragraph markup itself.\n\n$d = $word.ActiveDocument\n\n# --- Update the date paragraph ---\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd(\"`r`a\", \"`r\") -ne \"2023-07-27 Thursday\") {\n    throw \"Date paragraph text mismatch; expected 2023-07-27 Thursday\"\n}\n$dateParagraph.Range.Text = \"2023-07-28 Friday\"\n\n# --- Update every answer cell in the (single) table ---\n# Table.Cell(row, col) is 1-based, so add 1 to the 0-based coordinates.\n$table = $d.Tables.Item(1)\n\n$edits = @(\n    @{ Row = 1; Col = 1; Old = \"32\u00f74=8, 0\"; New = \"31\u00f75=6, 1\" },\n    @{ Row = 1; Col = 2; Old = \"29\u00f74=7, 1\"; New = \"11\u00f74=2, 3\" },\n    @{ Row = 1; Col = 3; Old = \"68\u00f79=7, 5\"; New = \"47\u00f78=5, 7\" },\n    @{ Row = 1; Col = 4; Old = \"12\u00f72=6, 0\"; New = \"48\u00f73=16, 0\" },\n    @{ Row = 1; Col = 5; Old = \"49\u00f73=16, 1\"; New = \"11\u00f72=5, 1\" },\n    @{ Row = 5; Col = 1; Old = \"46\u00f79=5, 1\"; New = \"82\u00f76=13, 4\" },\n    @{ Row = 5; Col = 2; Old = \"55\u00f77=7, 6\"; New = \"13\u00f78=1, 5\" },\n    @{ Row = 5; Col = 3; Old = \"68\u00f75=13, 3\"; New = \"48\u00f78=6, 0\" },\n    @{ Row = 5; Col = 4; Old = \"47\u00f78=5, 7\"; New = \"20\u00f78=2, 4\" },\n    @{ Row = 5; Col = 5; Old = \"59\u00f76=9, 5\"; New = \"80\u00f73=26, 2\" },\n    @{ Row = 9; Col = 1; Old = \"69\u00f73=23, 0\"; New = \"35\u00f78=4, 3\" },\n    @{ Row = 9; Col = 2; Old = \"50\u00f77=7, 1\"; New = \"14\u00f78=1, 6\" },\n    @{ Row = 9; Col = 3; Old = \"58\u00f78=7, 2\"; New = \"32\u00f76=5, 2\" },\n    @{ Row = 9; Col = 4; Old = \"31\u00f79=3, 4\"; New = \"82\u00f76=13, 4\" },\n    @{ Row = 9; Col = 5; Old = \"21\u00f77=3, 0\"; New = \"42\u00f79=4, 6\" },\n    @{ Row = 13; Col = 1; Old = \"52\u00f76=8, 4\"; New = \"84\u00f74=21, 0\" },\n    @{ Row = 13; Col = 2; Old = \"60\u00f76=10, 0\"; New = \"54\u00f78=6, 6\" },\n    @{ Row = 13; Col = 3; Old = \"31\u00f75=6, 1\"; New = \"99\u00f72=49, 1\" },\n    @{ Row = 13; Col = 4; Old = \"80\u00f73=26, 2\"; New = \"32\u00f73=10, 2\" },\n    @{ Row = 13; Col = 5; Old = \"82\u00f72=41, 0\"; New = \"98\u00f77=14, 0\" },\n    @{ Row = 17; Col = 1; Old = \"32\u00f77=4, 4\"; New = \"19\u00f74=4, 3\" },\n    @{ Row = 17; Col = 2; Old = \"81\u00f72=40, 1\"; New = \"92\u00f75=18, 2\" },\n    @{ Row = 17; Col = 3; Old = \"45\u00f78=5, 5\"; New = \"67\u00f75=13, 2\" },\n    @{ Row = 17; Col = 4; Old = \"94\u00f73=31, 1\"; New = \"75\u00f74=18, 3\" },\n    @{ Row = 17; Col = 5; Old = \"38\u00f72=19, 0\"; New = \"99\u00f78=12, 3\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $actual = $cell.Range.Text.TrimEnd(\"`r`a\", \"`r\")\n    if ($actual -ne $edit.Old) {\n        throw \"Cell ($($edit.Row), $($edit.Col)) text mismatch; expected $($edit.Old), found $actual\"\n    }\n    $cell.Range.Text = $edit.New\n}\n\n"}
